$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.719.59"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.647.96"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("E5").Value = "  +1.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.504"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.254"
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.43"
$ws.Range("E10").Value = "  +1.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0843"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.878.75"
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.677.44"
$ws.Range("E13").Value = "  +2.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.24"
$ws.Range("E14").Value = "  +3.49%  "
$ws.Range("E15").Value = "  +1.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.41"
$ws.Range("E16").Value = "  +5.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.772.01"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("E18").Value = "  +1.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "220.11"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.00"
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("E21").Value = "  +2.19%  "
$ws.Range("E22").Value = "  +2.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.59"
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.12"
$ws.Range("E24").Value = "  +9.94%  "
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.11"
$ws.Range("E28").Value = "  +2.84%  "
$ws.Range("E29").Value = "  +2.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0520"
$ws.Range("E30").Value = "  +1.88%  "
$ws.Range("E31").Value = "  +1.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.42"
$ws.Range("E32").Value = "  +3.39%  "
$ws.Range("E33").Value = "  +3.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.289.53"
$ws.Range("E34").Value = "  +7.45%  "
$ws.Range("E35").Value = "  +2.90%  "
$ws.Range("E36").Value = "  +7.20%  "
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.833"
$ws.Range("E38").Value = "  +2.86%  "
$ws.Range("E39").Value = "  +4.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.811"
$ws.Range("E41").Value = "  +2.26%  "
$ws.Range("E42").Value = "  -1.76%  "
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.789.45"
$ws.Range("E44").Value = "  +1.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.83"
$ws.Range("E45").Value = "  +1.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.94"
$ws.Range("E46").Value = "  +9.45%  "
$ws.Range("E47").Value = "  +3.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0517"
$ws.Range("E48").Value = "  +0.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.80"
$ws.Range("E49").Value = "  +2.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0979"
$ws.Range("E50").Value = "  +3.43%  "
